# Fall '22 week 9 "day after" update: add a new "Week 43" column of scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added week column (shared string "Week 43").
$ws.Range("AR1").Value = "Week 43"

# New scores reported for this week (players on rows 5, 7, 8, 9).
$ws.Range("AR5").Value = 10
$ws.Range("AR7").Value = 2
$ws.Range("AR8").Value = 3
$ws.Range("AR9").Value = 1.5

# Move the active selection to the newly-entered cell, matching the
# saved workbook view state after the edit.
$ws.Range("AR10").Select()
